$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# --- Clean up stray header cells C1:F1 (leftover "value" labels, no longer present) ---
$ws.Range("C1:F1").ClearContents()

# --- Insert a new row 9 for the "L_curve" parameter (pushes old rows 9-16 down to 10-17) ---
$ws.Rows.Item(9).Insert()

# Row 8 becomes "production_function" (was "Model"); value stays "Sigmoid"
$ws.Range("A1").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value2 = "production_function"

# Row 9 (new): "L_curve" label / numeric value 1 formatted like the other "value" column cells
$ws.Range("A1").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value2 = "L_curve"

$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value2 = 1

# --- Remove the old "Deletion" row (now shifted down to row 17) ---
$ws.Rows.Item(17).Delete()

# --- Update sheet view: selection on B10, make this the active sheet/tab ---
$ws.Range("B10").Select()
$ws.Activate()
